# Apply data corrections to the "Inscricoes" sheet.
# Each entry updates a single numeric cell (Inscritos/Pagos/Inscrições homologadas columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 52

$ws.Range("F4").Value = 6
$ws.Range("H4").Value = 20

$ws.Range("E17").Value = 124
$ws.Range("F17").Value = 60
$ws.Range("H17").Value = 92

$ws.Range("F18").Value = 49
$ws.Range("H18").Value = 85

$ws.Range("F24").Value = 15
$ws.Range("H24").Value = 19

$ws.Range("F26").Value = 16
$ws.Range("H26").Value = 26

$ws.Range("E33").Value = 43

$ws.Range("E36").Value = 105

$ws.Range("F42").Value = 18
$ws.Range("H42").Value = 27

$ws.Range("E43").Value = 25
$ws.Range("F43").Value = 15
$ws.Range("H43").Value = 18

$ws.Range("E48").Value = 32
$ws.Range("F48").Value = 22
$ws.Range("H48").Value = 27

$ws.Range("E49").Value = 73

$ws.Range("E51").Value = 12

$ws.Range("E62").Value = 44
$ws.Range("F62").Value = 11
$ws.Range("H62").Value = 25

$ws.Range("E63").Value = 37

$ws.Range("F64").Value = 18
$ws.Range("H64").Value = 23

$ws.Range("E65").Value = 31

$ws.Range("F66").Value = 22
$ws.Range("H66").Value = 30

$ws.Range("F67").Value = 23
$ws.Range("H67").Value = 32

$ws.Range("E71").Value = 37

$ws.Range("F72").Value = 22
$ws.Range("H72").Value = 33

$ws.Range("E76").Value = 53
$ws.Range("F76").Value = 20
$ws.Range("H76").Value = 37

$ws.Range("F77").Value = 23
$ws.Range("H77").Value = 40

$ws.Range("E80").Value = 29
$ws.Range("F80").Value = 11
$ws.Range("H80").Value = 23

$ws.Range("F81").Value = 11
$ws.Range("H81").Value = 16

$ws.Range("F88").Value = 12
$ws.Range("H88").Value = 20
